$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column Q data (customer_password_query) ---
$ws.Range("Q1").Value = "customer_password_query"
$ws.Range("Q1").NumberFormat = "@"
$ws.Range("Q2").Value = "Select count(*) from (Select MAX(AA.CREATED_ON) from DC_CUSTOMER_PASSWORD_HISTORY AA INNER JOIN DC_CUSTOMER_INFO BB ON AA.CUSTOMER_INFO_ID = BB.CUSTOMER_INFO_ID where AA.CUSTOMER_INFO_ID = (Select CUSTOMER_INFO_ID from dc_customer_info l where L.CUSTOMER_NAME = '{customer_name}') and AA.TRANSACTION_TYPE_ID = (Select LL.TRANSACTION_TYPE_ID from DC_TRANSACTION LL where LL.TRANSACTION_ID = '{TRANSACTION_ID}') and AA.PASSWORD = BB.TRANSACTION_PASSWORD and TRUNC(AA.CREATED_ON) < (SELECT TRUNC(SYSDATE) FROM DUAL) and TRUNC(AA.UPDATED_ON) < (SELECT TRUNC(SYSDATE) FROM DUAL) order by AA.UPDATED_ON desc) where rownum = 1"

# --- Size the new column to match the sheet's bestFit-style widths ---
$ws.Columns.Item(17).ColumnWidth = 24.43

# --- Update the view: scroll so column P is the left-most visible column
#     and select Q6 (matches the saved sheetView/selection in the target file) ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 16
$win.ScrollRow = 1
[void]$ws.Range("Q6").Select()
